$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("250 Waterdale Rd, Ivanhoe ...")
# shifting it (and everything below) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new exposure site.
$ws.Range("A8").Value = "13 Howard St, Altona Meadows VIC 3028"
$ws.Range("B8").Value = -37.878812
$ws.Range("C8").Value = 144.781325
$ws.Range("D8").Value = "Hobsons Bay (C)"
